# Generate Report for Handoff
# Updates the status of "9b57617c-67ae-43af-9945-b836536d0c0b.md" from
# "Handed back: in sync with en-US" to "Ready for handoff" on the
# Overview, zh-cn, and de-de sheets, along with the associated
# handoff datetime stamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 corresponds to 9b57617c-67ae-43af-9945-b836536d0c0b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-24 15:03:38"

# --- zh-cn sheet: row 3 corresponds to 9b57617c-67ae-43af-9945-b836536d0c0b.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-24 15:03:33"

# --- de-de sheet: row 3 corresponds to 9b57617c-67ae-43af-9945-b836536d0c0b.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-24 15:03:38"
